$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.572.17"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.034"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +3.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.029"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4382"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07405"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8762"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.862.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.516"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.697"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07205"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.036"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  +2.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.582.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.263"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.077.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.948"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.305"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.211"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7681"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.530"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.892"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.155"
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01979"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.825"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5177"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1673"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.737"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.593"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "108.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.713"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4664"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06393"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.889"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.26%  "
